$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column in H1, matching the style of existing headers (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill H2:H7 with 0 values (numeric), matching data rows
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
